# Applies the diff: insert a new data row before current row 27 (Vega Modelo de
# Temuco - Camote), which shifts all subsequent rows (27-85) down by one
# (becoming 28-86), and populates the new row 27 with its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27; Excel automatically shifts rows 27:85
# down to 28:86 and extends the sheet dimension to A1:R86.
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with its data.
$ws.Range("A27").Value() = 10
$ws.Range("B27").Value() = "Vega Modelo de Temuco"
$ws.Range("C27").Value() = "La Araucanía"
$ws.Range("D27").Value() = 44708
$ws.Range("E27").Value() = 9
$ws.Range("F27").Value() = 100114002
$ws.Range("G27").Value() = "Camote"
$ws.Range("H27").Value() = "Sin especificar"
$ws.Range("I27").Value() = "Primera"
$ws.Range("J27").Value() = 30
$ws.Range("K27").Value() = 18000
$ws.Range("L27").Value() = 18000
$ws.Range("M27").Value() = 18000
$ws.Range("N27").Value() = "$/caja 15 kilos granel"
$ws.Range("O27").Value() = "Perú"
$ws.Range("P27").Value() = 1200
$ws.Range("Q27").Value() = 15
$ws.Range("R27").Value() = "Hortaliza"
